$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue 'D2' '66.059.21'
Set-TextValue 'E2' '  +1.88%  '

Set-TextValue 'D3' '3.191.07'
Set-TextValue 'E3' '  +1.51%  '

Set-TextValue 'E4' '  -0.22%  '

Set-TextValue 'D5' '594.56'
Set-TextValue 'E5' '  +3.85%  '

Set-TextValue 'D6' '154.15'
Set-TextValue 'E6' '  +3.59%  '

Set-TextValue 'D7' '1.00'
Set-TextValue 'E7' '  -0.09%  '

Set-TextValue 'D8' '3.187.36'
Set-TextValue 'E8' '  +1.51%  '

Set-TextValue 'D9' '0.534'
Set-TextValue 'E9' '  +1.91%  '

Set-TextValue 'E10' '  +0.94%  '

Set-TextValue 'E11' '  +0.14%  '

Set-TextValue 'D12' '0.513'
Set-TextValue 'E12' '  +3.45%  '

Set-TextValue 'E13' '  +3.27%  '

Set-TextValue 'D14' '38.93'
Set-TextValue 'E14' '  +5.58%  '

Set-TextValue 'D15' '3.713.48'
Set-TextValue 'E15' '  +1.16%  '

Set-TextValue 'D16' '66.008.23'
Set-TextValue 'E16' '  +1.51%  '

Set-TextValue 'D17' '7.42'
Set-TextValue 'E17' '  +5.10%  '

Set-TextValue 'D18' '3.194.61'
Set-TextValue 'E18' '  +0.87%  '

Set-TextValue 'E19' '  +0.48%  '

Set-TextValue 'D20' '509.85'
Set-TextValue 'E20' '  +0.84%  '

Set-TextValue 'E21' '  +4.02%  '

Set-TextValue 'D22' '0.740'
Set-TextValue 'E22' '  +3.83%  '

Set-TextValue 'E23' '  -0.05%  '

Set-TextValue 'E24' '  +3.96%  '

Set-TextValue 'D25' '84.76'
Set-TextValue 'E25' '  +0.88%  '

Set-TextValue 'E26' '  +0.10%  '

Set-TextValue 'D27' '9.32'
Set-TextValue 'E27' '  +5.61%  '

Set-TextValue 'D28' '2.99'
Set-TextValue 'E28' '  +3.37%  '

Set-TextValue 'D29' '2.27'
Set-TextValue 'E29' '  +5.74%  '

Set-TextValue 'B30' 'Stacks'
Set-TextValue 'C30' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D30' '2.90'
Set-TextValue 'E30' '  +4.98%  '

Set-TextValue 'B31' 'NEARProtocol'
Set-TextValue 'C31' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D31' '6.92'
Set-TextValue 'E31' '  +11.94%  '

Set-TextValue 'D32' '28.25'
Set-TextValue 'E32' '  +2.93%  '

Set-TextValue 'D33' '1.22'
Set-TextValue 'E33' '  +3.65%  '

Set-TextValue 'D34' '0.999'
Set-TextValue 'E34' '  -0.17%  '

Set-TextValue 'D35' '6.55'
Set-TextValue 'E35' '  +1.14%  '

Set-TextValue 'E36' '  -0.63%  '

Set-TextValue 'D37' '0.0903'
Set-TextValue 'E37' '  +1.69%  '

Set-TextValue 'D38' '483.38'
Set-TextValue 'E38' '  +4.64%  '

Set-TextValue 'D39' '0.0417'
Set-TextValue 'E39' '  -0.64%  '

Set-TextValue 'D40' '2.92'
Set-TextValue 'E40' '  -1.60%  '

Set-TextValue 'D41' '8.83'
Set-TextValue 'E41' '  +2.58%  '

Set-TextValue 'D42' '0.122'
Set-TextValue 'E42' '  +4.65%  '

Set-TextValue 'D43' '0.297'
Set-TextValue 'E43' '  +6.05%  '

Set-TextValue 'D44' '0.0₃0651'
Set-TextValue 'E44' '  +12.55%  '

Set-TextValue 'D45' '2.929.79'
Set-TextValue 'E45' '  -4.05%  '

Set-TextValue 'D46' '2.42'
Set-TextValue 'E46' '  +0.72%  '

Set-TextValue 'D47' '28.46'
Set-TextValue 'E47' '  +0.24%  '

Set-TextValue 'E48' '  -0.03%  '

Set-TextValue 'E49' '  +2.37%  '

Set-TextValue 'E50' '  +4.23%  '

Set-TextValue 'E51' '  +7.34%  '
